# Append two new daily rows (2025-11-11 / serial 45972) for both charging
# stations to the bottom of the data table on sheet1, then move the active
# selection the way the workbook had it after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: 四方坪站 (station 1) for date serial 45972
$ws.Range("A22").Value = 45972
$ws.Range("B22").Value = "四方坪站"
$ws.Range("C22").Value = 8716.41
$ws.Range("D22").Value = 7621.81
$ws.Range("E22").Value = 2990.34
$ws.Range("F22").Value = 378

# Row 23: 高岭站 (station 2) for date serial 45972
$ws.Range("A23").Value = 45972
$ws.Range("B23").Value = "高岭站"
$ws.Range("C23").Value = 4125.0600000000004
$ws.Range("D23").Value = 3662.51
$ws.Range("E23").Value = 1140.6600000000001
$ws.Range("F23").Value = 166

# Reflect the final selected cell left in the worksheet after the edit.
[void]$ws.Range("I19").Select()
